# tieosoitteilla_ajanhetki.xlsx
# "Added _TULOS to result file name AND added Tunniste-field"
#
# The file-name change lives outside this workbook (in the exporting code),
# so the only in-workbook change is the new "Tunniste" column header on
# Sheet1 (column G), plus a couple of small view/format tweaks that came
# along with the edit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New header cell G1 = "Tunniste", formatted like the other header cells
# (A1:D1 use style index 1 - bold-ish header font). Copy the format from
# A1 so G1 ends up with the identical cell style, then set its text.
$ws1.Range("A1").Copy() | Out-Null
$ws1.Range("G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws1.Range("G1").Value = "Tunniste"

# Column F was widened slightly once the sheet had a 7th column.
$ws1.Columns.Item(6).ColumnWidth = 10.8

# Leave the selection where the editor last clicked.
$ws1.Range("J10").Select() | Out-Null

# Sheet2 / Sheet3 picked up an explicit print/page setup (paper size +
# portrait orientation) when the workbook was resaved in a newer Excel.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1
